$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin rows whose updated Price/Volume values differ from the scraped snapshot.
# Price (column D) values that would otherwise be auto-parsed as numbers are
# first switched to Text format so they keep being stored as strings (matching
# the original inline-string cells), then the style is reset to Normal so no
# stray number-format style is left behind.

$ws.Range("D2").Value = "70.950.44"
$ws.Range("E2").Value = "  +2.20%  "

$ws.Range("D3").Value = "3.588.79"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.18%  "

$ws.Range("D7").Value = "3.583.85"
$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  +4.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.592"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000280"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "4.167.64"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "621.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("D18").Value = "3.590.07"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").Value = "71.027.19"
$ws.Range("E19").Value = "  +2.24%  "

$ws.Range("E20").Value = "  -2.97%  "

$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  -16.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("E32").Value = "  -2.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.72%  "

$ws.Range("E34").Value = "  -2.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "630.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.34%  "

$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0488"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("E42").Value = "  +4.27%  "

$ws.Range("D43").Value = "3.420.02"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("D45").Value = "0.0₃0721"
$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("E51").Value = "  -0.02%  "
